$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

Set-TextValue "D2" "68.190.12"
Set-TextValue "E2" "  -2.08%  "
Set-TextValue "D3" "3.828.67"
Set-TextValue "E3" "  -1.48%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "600.84"
Set-TextValue "E5" "  -0.79%  "
Set-TextValue "D6" "169.58"
Set-TextValue "E6" "  -0.27%  "
Set-TextValue "D7" "3.825.89"
Set-TextValue "E7" "  -1.53%  "
Set-TextValue "E8" "  -0.04%  "
Set-TextValue "D9" "0.527"
Set-TextValue "E9" "  -1.55%  "
Set-TextValue "E10" "  -2.59%  "
Set-TextValue "D11" "6.46"
Set-TextValue "E11" "  +1.28%  "
Set-TextValue "D12" "0.457"
Set-TextValue "E12" "  -2.29%  "
Set-TextValue "E13" "  +3.15%  "
Set-TextValue "D14" "37.12"
Set-TextValue "E14" "  -3.01%  "
Set-TextValue "D15" "4.478.90"
Set-TextValue "E15" "  -1.38%  "
Set-TextValue "D16" "3.830.17"
Set-TextValue "E16" "  -1.86%  "
Set-TextValue "D17" "68.317.65"
Set-TextValue "E17" "  -1.83%  "
Set-TextValue "D18" "18.45"
Set-TextValue "E18" "  -1.51%  "
Set-TextValue "D19" "7.40"
Set-TextValue "E19" "  -3.31%  "
Set-TextValue "E20" "  -0.87%  "
Set-TextValue "D21" "11.08"
Set-TextValue "E21" "  -0.39%  "
Set-TextValue "D22" "468.40"
Set-TextValue "E22" "  -4.38%  "
Set-TextValue "D23" "0.737"
Set-TextValue "E23" "  -1.72%  "
Set-TextValue "D24" "0.0000160"
Set-TextValue "E24" "  -3.97%  "
Set-TextValue "D25" "83.13"
Set-TextValue "E25" "  -2.57%  "
Set-TextValue "D26" "2.23"
Set-TextValue "E26" "  -3.12%  "
Set-TextValue "D27" "12.22"
Set-TextValue "E27" "  -0.88%  "
Set-TextValue "D28" "10.03"
Set-TextValue "E28" "  -1.21%  "
Set-TextValue "E29" "  +0.00%  "
Set-TextValue "E30" "  -0.20%  "
Set-TextValue "D31" "3.977.79"
Set-TextValue "E31" "  -1.46%  "
Set-TextValue "D32" "7.69"
Set-TextValue "E32" "  -1.67%  "
Set-TextValue "D33" "31.58"
Set-TextValue "E33" "  -1.12%  "
Set-TextValue "D34" "2.31"
Set-TextValue "E34" "  -4.38%  "
Set-TextValue "D35" "9.41"
Set-TextValue "E35" "  -2.13%  "
Set-TextValue "D36" "3.792.92"
Set-TextValue "E36" "  -1.39%  "
Set-TextValue "E37" "  -2.06%  "
Set-TextValue "D38" "3.69"
Set-TextValue "E38" "  +10.32%  "
Set-TextValue "E39" "  -0.71%  "
Set-TextValue "E40" "  -1.96%  "
Set-TextValue "D41" "5.95"
Set-TextValue "E41" "  -3.09%  "
Set-TextValue "E42" "  +0.21%  "
Set-TextValue "D43" "0.315"
Set-TextValue "E43" "  -4.10%  "
Set-TextValue "D44" "1.98"
Set-TextValue "E44" "  -5.56%  "
Set-TextValue "D45" "8.75"
Set-TextValue "E45" "  +0.48%  "
Set-TextValue "D49" "47.14"
Set-TextValue "E49" "  -2.14%  "

# Rows 46-51 reorder + data refresh
Set-TextValue "B46" "Bittensor"
Set-TextValue "C46" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D46" "417.44"
Set-TextValue "E46" "  -4.47%  "
Set-TextValue "B47" "USDe"
Set-TextValue "C47" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue "D47" "1.00"
Set-TextValue "E47" "  -0.02%  "
Set-TextValue "B48" "FLOKI"
Set-TextValue "C48" "https://coinranking.com/coin/fmHk13Rqw+floki-floki"
Set-TextValue "D48" "0.000292"
Set-TextValue "E48" "  +6.42%  "
Set-TextValue "B50" "Monero"
Set-TextValue "C50" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D50" "141.87"
Set-TextValue "E50" "  -1.36%  "
Set-TextValue "B51" "VeChain"
Set-TextValue "C51" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D51" "0.0360"
Set-TextValue "E51" "  -2.72%  "
